$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 1440
$ws.Range("AB2").Value = 1487
$ws.Range("AF2").Value = 29186
$ws.Range("AA3").Value = 1643
$ws.Range("AB3").Value = 1715
$ws.Range("AF3").Value = 32410
$ws.Range("AA4").Value = 1571
$ws.Range("AB4").Value = 2116
$ws.Range("AF4").Value = 33956
$ws.Range("AA5").Value = 1567
$ws.Range("AB5").Value = 2295
$ws.Range("AF5").Value = 34766
$ws.Range("AA6").Value = 1919
$ws.Range("AB6").Value = 2656
$ws.Range("AF6").Value = 38181
$ws.Range("AA7").Value = 2171
$ws.Range("AB7").Value = 2950
$ws.Range("AF7").Value = 41055
$ws.Range("AA8").Value = 2295
$ws.Range("AB8").Value = 3284
$ws.Range("AF8").Value = 43687
$ws.Range("AA9").Value = 2216
$ws.Range("AB9").Value = 3662
$ws.Range("AF9").Value = 47903
$ws.Range("AA10").Value = 2480
$ws.Range("AB10").Value = 4404
$ws.Range("AF10").Value = 55506
$ws.Range("AA11").Value = 2646
$ws.Range("AB11").Value = 5256
$ws.Range("AF11").Value = 62990
$ws.Range("AA12").Value = 2971
$ws.Range("AB12").Value = 6390
$ws.Range("AF12").Value = 75673
$ws.Range("AA13").Value = 3909
$ws.Range("AB13").Value = 7477
$ws.Range("AF13").Value = 83577
$ws.Range("AA14").Value = 4682
$ws.Range("AB14").Value = 8586
$ws.Range("AF14").Value = 85894
$ws.Range("AA15").Value = 4981
$ws.Range("AB15").Value = 9599
$ws.Range("AF15").Value = 89082
$ws.Range("AA16").Value = 5223
$ws.Range("AB16").Value = 10292
$ws.Range("AF16").Value = 102688
$ws.Range("AA17").Value = 5925
$ws.Range("AB17").Value = 12040
$ws.Range("AF17").Value = 112075
$ws.Range("AA18").Value = 6120
$ws.Range("AB18").Value = 13705
$ws.Range("AF18").Value = 118873
